# "Added some smaller airlines"
#
# Fills in the missing airline Name (column C) for several smaller
# carriers further down the MAN flights table, and leaves behind the
# handful of hidden "_xlchart.v1.*" defined names that Excel creates as
# bookkeeping artifacts whenever the Insert/Recommended Charts gallery is
# invoked against the B1:F22 data block (even though no chart ends up
# being kept on the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Populate the previously-blank airline names in column C for the
#    smaller carriers (rows 23-36). Row 35 (EJU) intentionally has no
#    name filled in, matching the source data.
# ---------------------------------------------------------------------
$ws.Range("C23").Value = "Finnair"
$ws.Range("C24").Value = "Norwegian Air Shuttle"
$ws.Range("C25").Value = "Singapore Airlines"
$ws.Range("C26").Value = "Pakistan International Airlines"
$ws.Range("C27").Value = "Icelandair"
$ws.Range("C28").Value = "Oman Air"
$ws.Range("C29").Value = "Pegasus Airlines"
$ws.Range("C30").Value = "TAP Portugal"
$ws.Range("C31").Value = "United Airlines"
$ws.Range("C32").Value = "Vueling Airlines"
$ws.Range("C33").Value = "American Airlines"
$ws.Range("C34").Value = "Austrian Airlines"
$ws.Range("C36").Value = "Federal Express"

# ---------------------------------------------------------------------
# 2. Recreate the hidden "_xlchart" defined names left behind by the
#    Quick Analysis / Recommended Charts feature, which was invoked
#    three times over the B1:F22 block.
# ---------------------------------------------------------------------
$chartRefs = @(
  "Sheet1!`$B`$1",
  "Sheet1!`$B`$2:`$B`$22",
  "Sheet1!`$C`$1",
  "Sheet1!`$C`$2:`$C`$22",
  "Sheet1!`$D`$1",
  "Sheet1!`$D`$2:`$D`$22",
  "Sheet1!`$E`$1",
  "Sheet1!`$E`$2:`$E`$22",
  "Sheet1!`$F`$1",
  "Sheet1!`$F`$2:`$F`$22"
)

for ($batch = 0; $batch -lt 3; $batch++) {
    for ($i = 0; $i -lt $chartRefs.Length; $i++) {
        $idx = $batch * 10 + $i
        $defName = $wb.Names.Add("_xlchart.v1.$idx", "=" + $chartRefs[$i])
        $defName.Visible = $false
    }
}

# ---------------------------------------------------------------------
# 3. Leave the active selection on C34 (scrolled back to the top),
#    matching where the author ended up after editing the column.
# ---------------------------------------------------------------------
$ws.Range("C34").Select()
